$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simulate typing new job numbers, then overwriting them again,
# so that the shared string table accumulates the now-unused values
# while the cells end up referencing the final values.
$ws.Range("B2").Value = "32373945"
$ws.Range("B3").Value = "32373946"
$ws.Range("B4").Value = "32373947"

$ws.Range("B2").Value = "32376214"
$ws.Range("B3").Value = "32376215"
$ws.Range("B4").Value = "32376217"
